# Update the GitHub Link slide's body text from the README reminder
# to the actual repository URL.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "https://github.com/sanketrshinde19/-Steganography.git"
